$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.478.72"
$ws.Range("D3").Value = "1.821.60"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'311.70"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.4237"
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").Value = "'0.3623"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").Value = "'0.8567"
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("D11").Value = "'20.56"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").Value = "1.835.91"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "'5.390"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "'6.460"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "'0.06913"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").Value = "'1.005"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'80.09"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").Value = "'0.000008859"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "27.183.89"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Value = "'5.108"
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("D23").Value = "'10.88"
$ws.Range("E23").Value = "  +5.41%  "
$ws.Range("D24").Value = "2.025.68"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "'1.984"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").Value = "'18.68"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "'5.134"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").Value = "'113.88"
$ws.Range("E29").Value = "  -4.58%  "
$ws.Range("D30").Value = "'1.798"
$ws.Range("E30").Value = "  -3.65%  "
$ws.Range("D31").Value = "'0.08833"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").Value = "'2.983"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").Value = "'0.7411"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").Value = "'4.517"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "'1.119"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "'1.086"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "'0.05266"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").Value = "'0.01922"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'2.784"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D42").Value = "'0.1639"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("D43").Value = "'6.448"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "'10.35"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").Value = "'105.36"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Value = "'0.06448"
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("D48").Value = "'0.4651"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").Value = "'1.002"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "'1.608"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").Value = "'63.19"
$ws.Range("E51").Value = "  -1.58%  "
